# Applies the "Final testing data added and processed, OpenSim plots updated"
# commit: fills in the previously-blank measurement columns (B:E) for the
# STATIC and DYNAMIC sample rows (4, 5, 10, 11, 16, 17) which frees the
# shared F/G formulas from their #DIV/0! error state, and moves the active
# selection to F17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 4 (Sample 2 / STATIC block) ---
$ws.Range("B4").Value = 101
$ws.Range("C4").Value = 17.9
$ws.Range("D4").Value = 4.01
$ws.Range("E4").Value = 1.2

# --- Row 5 (Sample 3 / STATIC block) ---
$ws.Range("B5").Value = 72
$ws.Range("C5").Value = 17.97
$ws.Range("D5").Value = 3.34
$ws.Range("E5").Value = 1.23

# --- Row 10 (Sample 2 / CONTROL-adjacent block) ---
$ws.Range("B10").Value = 72
$ws.Range("C10").Value = 19.33
$ws.Range("D10").Value = 3.75
$ws.Range("E10").Value = 1.2

# --- Row 11 (Sample 3) ---
$ws.Range("B11").Value = 76
$ws.Range("C11").Value = 18.39
$ws.Range("D11").Value = 3.56
$ws.Range("E11").Value = 0.73

# --- Row 16 (Sample 2 / DYNAMIC block) ---
$ws.Range("B16").Value = 121
$ws.Range("C16").Value = 18.74
$ws.Range("D16").Value = 3.97
$ws.Range("E16").Value = 1.51

# --- Row 17 (Sample 3 / DYNAMIC block) ---
$ws.Range("B17").Value = 138
$ws.Range("C17").Value = 18.55
$ws.Range("D17").Value = 4.3
$ws.Range("E17").Value = 1.73

# Recalculate so the shared F/G formulas pick up real values instead of
# the cached #DIV/0! results.
$excel.CalculateFull()

# Move the active selection to F17, matching the author's final cursor
# position after reviewing the newly computed column.
[void]$ws.Range("F17").Select()
